# Commit: "update file with jgit"
# Cell E8 on the "Rules" sheet changes from "Good Morning" to "GIT UPDATE",
# and that cell becomes the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
